$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("G18").Value = 2.57
$ws.Range("H18").Value = 3.4
$ws.Range("I18").Value = 2.52
$ws.Range("J18").Value = 3.15
$ws.Range("K18").Value = 2.12
$ws.Range("Q18").Value = 2.05
$ws.Range("R18").Value = 1.72
$ws.Range("S18").Value = 1.4
$ws.Range("T18").Value = 2.72
$ws.Range("V18").Value = 1.88
$ws.Range("Y18").Value = 10
$ws.Range("AB18").Value = 35
$ws.Range("AJ18").Value = 26
$ws.Range("AP18").Value = 22
$ws.Range("AQ18").Value = 60
$ws.Range("AR18").Value = 100
$ws.Range("AS18").Value = 300
$ws.Range("AT18").Value = 2.72
$ws.Range("AU18").Value = 7.3
$ws.Range("AW18").Value = 4.4
$ws.Range("AY18").Value = 22
$ws.Range("BA18").Value = 100

# Row 27
$ws.Range("G27").Value = 2.2
$ws.Range("H27").Value = 3.2
$ws.Range("L27").Value = 3.5
$ws.Range("N27").Value = 7.1
$ws.Range("V27").Value = 1.95
$ws.Range("W27").Value = 7.5
$ws.Range("X27").Value = 10.5
$ws.Range("AC27").Value = 7.1
$ws.Range("AD27").Value = 6.3
$ws.Range("AE27").Value = 14
$ws.Range("AG27").Value = 9.5
$ws.Range("AJ27").Value = 37
$ws.Range("AN27").Value = 4.15
$ws.Range("AP27").Value = 19.5
$ws.Range("AR27").Value = 80
$ws.Range("AU27").Value = 6.9
$ws.Range("AX27").Value = 16
$ws.Range("BA27").Value = 100
$ws.Range("BB27").Value = 300

# Row 36
$ws.Range("G36").Value = 2.92
$ws.Range("H36").Value = 3.05
$ws.Range("I36").Value = 2.37
$ws.Range("J36").Value = 3.55
$ws.Range("K36").Value = 2
$ws.Range("L36").Value = 2.95
$ws.Range("O36").Value = 1.29
$ws.Range("P36").Value = 3
$ws.Range("Q36").Value = 1.9
$ws.Range("R36").Value = 1.8
$ws.Range("S36").Value = 1.42
$ws.Range("T36").Value = 2.47
$ws.Range("U36").Value = 1.62
$ws.Range("V36").Value = 2.02
$ws.Range("W36").Value = 9
$ws.Range("X36").Value = 15.5
$ws.Range("Z36").Value = 37
$ws.Range("AA36").Value = 26
$ws.Range("AB36").Value = 32
$ws.Range("AC36").Value = 9.25
$ws.Range("AD36").Value = 6
$ws.Range("AF36").Value = 55
$ws.Range("AG36").Value = 8.5
$ws.Range("AH36").Value = 12.5
$ws.Range("AJ36").Value = 26
$ws.Range("AK36").Value = 18.5
$ws.Range("AL36").Value = 26
$ws.Range("AM36").Value = 400
$ws.Range("AO36").Value = 16.5
$ws.Range("AP36").Value = 24
$ws.Range("AQ36").Value = 80
$ws.Range("AR36").Value = 120
$ws.Range("AS36").Value = 300
$ws.Range("AT36").Value = 2.45
$ws.Range("AV36").Value = 60
$ws.Range("AW36").Value = 4.3
$ws.Range("AX36").Value = 12.5
$ws.Range("AY36").Value = 19.5
$ws.Range("AZ36").Value = 50
$ws.Range("BA36").Value = 80
$ws.Range("BB36").Value = 250

# Row 37
$ws.Range("G37").Value = 2.27
$ws.Range("H37").Value = 2.95
$ws.Range("I37").Value = 3.2
$ws.Range("J37").Value = 2.85
$ws.Range("K37").Value = 1.98
$ws.Range("L37").Value = 3.8
$ws.Range("M37").Value = 1.03
$ws.Range("N37").Value = 6.95
$ws.Range("O37").Value = 1.39
$ws.Range("P37").Value = 2.57
$ws.Range("Q37").Value = 2.12
$ws.Range("R37").Value = 1.57
$ws.Range("S37").Value = 1.44
$ws.Range("T37").Value = 2.4
$ws.Range("U37").Value = 1.82
$ws.Range("V37").Value = 1.78
$ws.Range("W37").Value = 6.9
$ws.Range("X37").Value = 10.75
$ws.Range("Y37").Value = 9
$ws.Range("Z37").Value = 24
$ws.Range("AA37").Value = 20
$ws.Range("AB37").Value = 32
$ws.Range("AC37").Value = 7.4
$ws.Range("AD37").Value = 5.8
$ws.Range("AE37").Value = 15
$ws.Range("AF37").Value = 80
$ws.Range("AG37").Value = 8
$ws.Range("AH37").Value = 15.5
$ws.Range("AI37").Value = 11.5
$ws.Range("AJ37").Value = 45
$ws.Range("AK37").Value = 32
$ws.Range("AL37").Value = 45
$ws.Range("AM37").Value = 700
$ws.Range("AN37").Value = 4.05
$ws.Range("AO37").Value = 11.75
$ws.Range("AP37").Value = 20
$ws.Range("AQ37").Value = 50
$ws.Range("AR37").Value = 80
$ws.Range("AS37").Value = 250
$ws.Range("AT37").Value = 2.37
$ws.Range("AU37").Value = 7
$ws.Range("AV37").Value = 65
$ws.Range("AW37").Value = 5
$ws.Range("AX37").Value = 18.5
$ws.Range("AY37").Value = 27
$ws.Range("AZ37").Value = 100
$ws.Range("BA37").Value = 150
$ws.Range("BB37").Value = 350

# Row 42
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = 3.3
$ws.Range("I42").Value = 3.4
$ws.Range("J42").Value = 2.65
$ws.Range("L42").Value = 3.9
$ws.Range("O42").Value = 1.3
$ws.Range("P42").Value = 3.2
$ws.Range("Q42").Value = 1.91
$ws.Range("R42").Value = 1.83
$ws.Range("W42").Value = 7.3
$ws.Range("X42").Value = 9.5
$ws.Range("Z42").Value = 18.5
$ws.Range("AG42").Value = 10.5
$ws.Range("AH42").Value = 19
$ws.Range("AI42").Value = 11.75
$ws.Range("AJ42").Value = 50
$ws.Range("AL42").Value = 35
$ws.Range("AN42").Value = 3.95
$ws.Range("AO42").Value = 10.5
$ws.Range("AU42").Value = 7.1
$ws.Range("AW42").Value = 5.3
$ws.Range("AX42").Value = 19
$ws.Range("AY42").Value = 25
$ws.Range("BA42").Value = 120
$ws.Range("BB42").Value = 300
